$d = $word.ActiveDocument

# "Versi" + "on" -> replace-in-place with identical text "Version"; Word
# merges the two runs that get touched into a single run.
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version", 2)

# The trailing "." lives in its own run just after the (empty) _GoBack
# bookmark. Find it and delete it first so the bookmark isn't swallowed by
# a replace that spans across it.
$dot = $d.Content
$dot.Find.Execute(".", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$dot.Delete()

# Now grow "2" (in the run right before the bookmark) into "1." - this
# merges into a single " 1." run, with the bookmark left intact right
# after it, matching "Version 1." with the bookmark at the very end.
$d.Content.Find.Execute("2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1.", 2)
